# Weekly update: insert two new price rows (date 44505) at the top of the
# "Camote" price block (rows 615-616), pushing the existing data down by two
# rows. This matches the new xlsx dimension A1:R725 (was A1:R723).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 615 (shifts everything from old row 615 onward
# down by two rows: old 615 -> new 617, ..., old 723 -> new 725).
$ws.Rows("615:616").Insert()

# --- New row 615 ---
$ws.Range("A615").Value = 9
$ws.Range("B615").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C615").Value = "Metropolitana"
$ws.Range("D615").Value = 44505
$ws.Range("E615").Value = 13
$ws.Range("F615").Value = 100112045
$ws.Range("G615").Value = "Zapallo"
$ws.Range("H615").Value = "Camote"
$ws.Range("I615").Value = "1a nueva(o)"
$ws.Range("J615").Value = 160
$ws.Range("K615").Value = 700
$ws.Range("L615").Value = 800
$ws.Range("M615").Value = 750
$ws.Range("N615").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O615").Value = "Perú"
$ws.Range("P615").Value = 750
$ws.Range("Q615").Value = 1
$ws.Range("R615").Value = "Hortaliza"

# --- New row 616 ---
$ws.Range("A616").Value = 9
$ws.Range("B616").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C616").Value = "Metropolitana"
$ws.Range("D616").Value = 44505
$ws.Range("E616").Value = 13
$ws.Range("F616").Value = 100112045
$ws.Range("G616").Value = "Zapallo"
$ws.Range("H616").Value = "Camote"
$ws.Range("I616").Value = "2a nueva(o)"
$ws.Range("J616").Value = 97
$ws.Range("K616").Value = 500
$ws.Range("L616").Value = 500
$ws.Range("M616").Value = 500
$ws.Range("N616").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O616").Value = "Perú"
$ws.Range("P616").Value = 500
$ws.Range("Q616").Value = 1
$ws.Range("R616").Value = "Hortaliza"

Write-Output "rows inserted and populated"
